$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear fromUser on row 3 (B3) but keep the cell record present
$ws.Range("B3").Value = ""
$ws.Range("B3").Style = "Normal"

# Row 4: new REGISTERPROJECT request from YCHERN to ASFLI, PENDING, projectID 2
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "YCHERN"
$ws.Range("C4").Value = "ASFLI"
$ws.Range("D4").Value = "REGISTERPROJECT"
$ws.Range("E4").Value = "PENDING"
$ws.Range("F4").Value = 2
$ws.Range("G4").Style = "Normal"
$ws.Range("H4").Style = "Normal"

# Row 5: CHANGETITLE request from YCHERN to ASMADHUKUMAR, PENDING, projectID 1, newTitle "tes"
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "YCHERN"
$ws.Range("C5").Value = "ASMADHUKUMAR"
$ws.Range("D5").Value = "CHANGETITLE"
$ws.Range("E5").Value = "PENDING"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "tes"
$ws.Range("H5").Style = "Normal"

# Row 6: CHANGETITLE request from YCHERN to ASMADHUKUMAR, PENDING, projectID 1, newTitle "test"
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "YCHERN"
$ws.Range("C6").Value = "ASMADHUKUMAR"
$ws.Range("D6").Value = "CHANGETITLE"
$ws.Range("E6").Value = "PENDING"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "test"
$ws.Range("H6").Style = "Normal"
